$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column width formatting for new columns T (20) and U (21)
$ws.Range("D1:U1").ColumnWidth = 8.7109375

# Row 3: new bottom-border-only cells T3, U3 (match existing style of Q3:S3)
$ws.Range("T3:U3").Value = $null
$ws.Range("T3:U3").Style = $ws.Range("S3").Style

# Row 4: new year headers T4=2023, U4=2024
$ws.Range("T4").Value = 2023
$ws.Range("U4").Value = 2024
$ws.Range("T4:U4").Style = $ws.Range("S4").Style

# Row 5: new data values T5=10.8, U5=6.5
$ws.Range("T5").Value = 10.8
$ws.Range("U5").Value = 6.5
$ws.Range("T5:U5").Style = $ws.Range("S5").Style

# Row 5 height change
$ws.Range("A5").RowHeight = 41.25

# Clear selection back to A1 (the diff removes the stored selection of T5)
$ws.Range("A1").Select()
